$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("EVData")
$ws5.Range("B1").Value = "BatteryCapacity(kWh)"
$r = $ws5.Columns.Item(2).EntireColumn.AutoFit()
Write-Output "result: $r"
